$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current last data row (row 375).
# This pushes the old row 375 down to row 378 and creates three blank
# rows at 375-377 that will be filled with the pre-edit data that used
# to live in rows 372-374.
$ws.Rows.Item(375).Insert()
$ws.Rows.Item(375).Insert()
$ws.Rows.Item(375).Insert()

# --- New row 375: old row 372 data (before this week's update) ---
$ws.Range("A375").Value = 3
$ws.Range("B375").Value = "Femacal de La Calera"
$ws.Range("C375").Value = "Coquimbo"
$ws.Range("D375").Value = 44335
$ws.Range("E375").Value = 5
$ws.Range("F375").Value = 100112021
$ws.Range("G375").Value = "Ají"
$ws.Range("H375").Value = "Inferno"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 35
$ws.Range("K375").Value = 30000
$ws.Range("L375").Value = 30000
$ws.Range("M375").Value = 30000
$ws.Range("N375").Value = "$/caja 15 kilos"
$ws.Range("O375").Value = "Región de Arica y Parinacota"
$ws.Range("P375").Value = 2000
$ws.Range("Q375").Value = 15
$ws.Range("R375").Value = "Hortaliza"

# --- New row 376: old row 373 data (before this week's update) ---
$ws.Range("A376").Value = 3
$ws.Range("B376").Value = "Femacal de La Calera"
$ws.Range("C376").Value = "Coquimbo"
$ws.Range("D376").Value = 44552
$ws.Range("E376").Value = 5
$ws.Range("F376").Value = 100112021
$ws.Range("G376").Value = "Ají"
$ws.Range("H376").Value = "Americana (o)"
$ws.Range("I376").Value = "Primera"
$ws.Range("J376").Value = 73
$ws.Range("K376").Value = 22000
$ws.Range("L376").Value = 23000
$ws.Range("M376").Value = 22521
$ws.Range("N376").Value = "$/caja 15 kilos"
$ws.Range("O376").Value = "Limache"
$ws.Range("P376").Value = 1501
$ws.Range("Q376").Value = 15
$ws.Range("R376").Value = "Hortaliza"

# --- New row 377: old row 374 data (before this week's update) ---
$ws.Range("A377").Value = 3
$ws.Range("B377").Value = "Femacal de La Calera"
$ws.Range("C377").Value = "Coquimbo"
$ws.Range("D377").Value = 44544
$ws.Range("E377").Value = 5
$ws.Range("F377").Value = 100112021
$ws.Range("G377").Value = "Ají"
$ws.Range("H377").Value = "Americana (o)"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 73
$ws.Range("K377").Value = 21000
$ws.Range("L377").Value = 22000
$ws.Range("M377").Value = 21521
$ws.Range("N377").Value = "$/caja 15 kilos"
$ws.Range("O377").Value = "Limache"
$ws.Range("P377").Value = 1435
$ws.Range("Q377").Value = 15
$ws.Range("R377").Value = "Hortaliza"

# --- Row 372: refresh with the new week's data ---
$ws.Range("D372").Value = 44595
$ws.Range("H372").Value = "Americana (o)"
$ws.Range("J372").Value = 80
$ws.Range("K372").Value = 13500
$ws.Range("L372").Value = 14000
$ws.Range("M372").Value = 13750
$ws.Range("O372").Value = "Limache"
$ws.Range("P372").Value = 917

# --- Row 373: refresh with the new week's data ---
$ws.Range("D373").Value = 44595
$ws.Range("J373").Value = 35
$ws.Range("K373").Value = 23000
$ws.Range("M373").Value = 23000
$ws.Range("N373").Value = "$/caja 25 kilos"
$ws.Range("O373").Value = "Provincia de Limarí"
$ws.Range("P373").Value = 920
$ws.Range("Q373").Value = 25

# --- Row 374: refresh with the new week's data ---
$ws.Range("D374").Value = 44595
$ws.Range("I374").Value = "Segunda"
$ws.Range("J374").Value = 30
$ws.Range("K374").Value = 19000
$ws.Range("L374").Value = 19000
$ws.Range("M374").Value = 19000
$ws.Range("N374").Value = "$/caja 25 kilos"
$ws.Range("O374").Value = "Provincia de Limarí"
$ws.Range("P374").Value = 760
$ws.Range("Q374").Value = 25
